# Add two new columns, I ("I0") and J ("IF"), to the right of the
# existing "IP" column (H) on the single worksheet, filling in the
# header labels and the per-row numeric data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the look of the existing header row (bold text, thin border,
# centered alignment) by copying the formatting from the neighboring
# "IP" header cell (H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows (2-34) ---------------------------------------------------
$iValues = @(9,8,5,7,5,8,4,8,6,8,8,6,8,8,4,5,7,6,5,6,5,8,5,9,5,9,5,6,8,8,9,7,8)
$jValues = @(9,8,5,8,5,9,5,9,6,8,8,7,8,8,6,7,7,6,7,7,6,8,6,9,6,9,6,6,8,8,9,7,8)

for ($row = 2; $row -le 34; $row++) {
    $idx = $row - 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]   # column I
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]  # column J
}
